# Refresh crafting-leve profit data (columns H:N) pulled from the market-board price feed.
# These are plain cached values (no formulas) updated row-by-row per the latest snapshot;
# rows whose HQ/NQ price no longer applies have that profit cell cleared entirely rather
# than zeroed, matching how the source sheet represents "not applicable".

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H10").Value = 15500
$ws.Range("J10").Value = 15500
$ws.Range("L10").Value = 15500
$ws.Range("N10").Value = -16086
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H20").Value = 1305.3334
$ws.Range("I20").Value = 566.4
$ws.Range("K20").Value = 566.4
$ws.Range("M20").Value = -336.4
$ws.Range("H35").Value = 1305.3334
$ws.Range("I35").Value = 566.4
$ws.Range("K35").Value = 566.4
$ws.Range("M35").Value = -187.4
$ws.Range("H107").Value = 714.2308
$ws.Range("I107").Value = 573
$ws.Range("K107").Value = 573
$ws.Range("M107").Value = 1347
$ws.Range("H116").Value = 6279.8
$ws.Range("I116").Value = 5726.2666
$ws.Range("J116").Value = 6833.3335
$ws.Range("K116").Value = 5726.2666
$ws.Range("L116").Value = 6833.3335
$ws.Range("M116").Value = -2284.2666
$ws.Range("N116").Value = -13717.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 3336666.2
$ws.Range("I10").Value = 3336666.2
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 3336666.2
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -3336496.2
$ws.Range("N10").ClearContents()
$ws.Range("H21").Value = 2167.1667
$ws.Range("I21").Value = 2051
$ws.Range("K21").Value = 2051
$ws.Range("M21").Value = -1677
$ws.Range("H36").Value = 13598.333
$ws.Range("I36").Value = 13598.333
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 13598.333
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -13252.333
$ws.Range("N36").ClearContents()
$ws.Range("H45").Value = 1713.7142
$ws.Range("I45").Value = 1666
$ws.Range("K45").Value = 1666
$ws.Range("M45").Value = -1289
$ws.Range("H132").Value = 3764.5334
$ws.Range("I132").Value = 2757.8572
$ws.Range("J132").Value = 4645.375
$ws.Range("K132").Value = 8273.571599999999
$ws.Range("L132").Value = 13936.125
$ws.Range("M132").Value = -5743.571599999999
$ws.Range("N132").Value = -18996.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2349.9167
$ws.Range("I134").Value = 2095.111
$ws.Range("K134").Value = 6285.333
$ws.Range("M134").Value = -3750.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 359.85715
$ws.Range("I5").Value = 176.5
$ws.Range("J5").Value = 433.2
$ws.Range("K5").Value = 176.5
$ws.Range("L5").Value = 433.2
$ws.Range("M5").Value = -64.5
$ws.Range("N5").Value = -657.2
$ws.Range("H6").Value = 1250374.8
$ws.Range("I6").Value = 1666999.6
$ws.Range("K6").Value = 1666999.6
$ws.Range("M6").Value = -1666886.6
$ws.Range("H22").Value = 1635.1666
$ws.Range("I22").Value = 1173.5
$ws.Range("K22").Value = 1173.5
$ws.Range("M22").Value = -823.5
$ws.Range("H25").Value = 1417.75
$ws.Range("I25").Value = 1257
$ws.Range("K25").Value = 1257
$ws.Range("M25").Value = -1083
$ws.Range("H31").Value = 1823.3182
$ws.Range("I31").Value = 1641.5555
$ws.Range("J31").Value = 2641.25
$ws.Range("K31").Value = 1641.5555
$ws.Range("L31").Value = 2641.25
$ws.Range("M31").Value = -1346.5555
$ws.Range("N31").Value = -3231.25
$ws.Range("H34").Value = 1823.3182
$ws.Range("I34").Value = 1641.5555
$ws.Range("J34").Value = 2641.25
$ws.Range("K34").Value = 1641.5555
$ws.Range("L34").Value = 2641.25
$ws.Range("M34").Value = -1439.5555
$ws.Range("N34").Value = -3045.25
$ws.Range("H36").Value = 2500
$ws.Range("J36").Value = 2500
$ws.Range("L36").Value = 2500
$ws.Range("N36").Value = -3276
$ws.Range("H40").Value = 2500
$ws.Range("J40").Value = 2500
$ws.Range("L40").Value = 2500
$ws.Range("N40").Value = -2820
$ws.Range("H99").Value = 7608.7
$ws.Range("I99").Value = 6865.857
$ws.Range("J99").Value = 9342
$ws.Range("K99").Value = 6865.857
$ws.Range("L99").Value = 9342
$ws.Range("M99").Value = -5367.857
$ws.Range("N99").Value = -12338
$ws.Range("H126").Value = 7608.7
$ws.Range("I126").Value = 6865.857
$ws.Range("J126").Value = 9342
$ws.Range("K126").Value = 20597.571
$ws.Range("L126").Value = 28026
$ws.Range("M126").Value = -18127.571
$ws.Range("N126").Value = -32966
$ws.Range("H138").Value = 10785
$ws.Range("I138").Value = 1499
$ws.Range("J138").Value = 27499.8
$ws.Range("K138").Value = 1499
$ws.Range("L138").Value = 27499.8
$ws.Range("M138").Value = 3641
$ws.Range("N138").Value = -37779.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 412
$ws.Range("I16").Value = 412.5
$ws.Range("J16").Value = 410
$ws.Range("K16").Value = 1237.5
$ws.Range("L16").Value = 1230
$ws.Range("M16").Value = -1064.5
$ws.Range("N16").Value = -1576
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H131").Value = 1638.5385
$ws.Range("I131").Value = 974.8570999999999
$ws.Range("J131").Value = 1883.0526
$ws.Range("K131").Value = 2924.5713
$ws.Range("L131").Value = 5649.1578
$ws.Range("M131").Value = 2115.4287
$ws.Range("N131").Value = -15729.1578

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2694.1667
$ws.Range("I80").Value = 2504.8572
$ws.Range("J80").Value = 2959.2
$ws.Range("K80").Value = 2504.8572
$ws.Range("L80").Value = 2959.2
$ws.Range("M80").Value = -1506.8572
$ws.Range("N80").Value = -4955.2
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").ClearContents()
$ws.Range("H83").Value = 2694.1667
$ws.Range("I83").Value = 2504.8572
$ws.Range("J83").Value = 2959.2
$ws.Range("K83").Value = 12524.286
$ws.Range("L83").Value = 14796
$ws.Range("M83").Value = -7532.286
$ws.Range("N83").Value = -24780
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 586.1667
$ws.Range("I9").Value = 279.75
$ws.Range("J9").Value = 1199
$ws.Range("K9").Value = 279.75
$ws.Range("L9").Value = 1199
$ws.Range("M9").Value = -55.75
$ws.Range("N9").Value = -1647
$ws.Range("H25").Value = 30500
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 30500
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 30500
$ws.Range("N25").Value = -30960
$ws.Range("M25").ClearContents()
$ws.Range("H26").Value = 3336.75
$ws.Range("I26").Value = 3099
$ws.Range("J26").Value = 3370.7144
$ws.Range("K26").Value = 3099
$ws.Range("L26").Value = 3370.7144
$ws.Range("M26").Value = -2804
$ws.Range("N26").Value = -3960.7144
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()
$ws.Range("H30").Value = 446.625
$ws.Range("I30").Value = 446.625
$ws.Range("K30").Value = 446.625
$ws.Range("M30").Value = -338.625
$ws.Range("H61").Value = 1750.7693
$ws.Range("I61").Value = 1374.5
$ws.Range("K61").Value = 1374.5
$ws.Range("M61").Value = -1172.5
$ws.Range("H113").Value = 1750.7693
$ws.Range("I113").Value = 1374.5
$ws.Range("K113").Value = 1374.5
$ws.Range("M113").Value = 795.5
$ws.Range("H136").Value = 4846.3335
$ws.Range("I136").Value = 4625.6
$ws.Range("K136").Value = 13876.8
$ws.Range("M136").Value = -11326.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 45000000
$ws.Range("J3").Value = 45000000
$ws.Range("L3").Value = 45000000
$ws.Range("N3").Value = -45000228
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H107").Value = 542.38464
$ws.Range("J107").Value = 589
$ws.Range("L107").Value = 1767
$ws.Range("N107").Value = -5607
$ws.Range("H136").Value = 5000
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
